$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update title
$ws.Range("A1").Value = "Model: forecasting_prophet_splited_model"

# Update header row (row 3)
$ws.Range("B3").Value = "mae"
$ws.Range("C3").Value = "mse"
$ws.Range("D3").Value = "rmse"
$ws.Range("E3").Value = "mape"

# Clear old extra header cells (F3:G3) and metric rows below
$ws.Range("F3:G9").Clear()
$ws.Range("A5:E9").Clear()

# Update data row (row 4)
$ws.Range("A4").Value = "metrics"
$ws.Range("B4").Value = 7612.847821239554
$ws.Range("C4").Value = 81425740.7853
$ws.Range("D4").Value = 9023.621267833663
$ws.Range("E4").Value = 25.6161856113109

# Clear old F4/G4 cells (roc_auc/threshold values) since columns no longer needed
$ws.Range("F4:G4").Clear()
